$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple frequency count bumps (+1)
$ws.Range("C8").Value = 1544
$ws.Range("C22").Value = 777
$ws.Range("C25").Value = 742
$ws.Range("C38").Value = 525
$ws.Range("C48").Value = 430
$ws.Range("C55").Value = 390

# Rows 110-113: words shift down one slot, with new word "man" inserted at top
$ws.Range("B110").Value = "man"
$ws.Range("C110").Value = 273

$ws.Range("B111").Value = "communist"
$ws.Range("C111").Value = 272

$ws.Range("B112").Value = "punch"
$ws.Range("C112").Value = 271

$ws.Range("B113").Value = "noth"
$ws.Range("C113").Value = 271

# Rows 119-122: words shift down one slot, with new word "read" inserted at top
$ws.Range("B119").Value = "read"
$ws.Range("C119").Value = 264

$ws.Range("B120").Value = "hitler"
$ws.Range("C120").Value = 264

$ws.Range("B121").Value = "issu"
$ws.Range("C121").Value = 263

$ws.Range("B122").Value = "stop"
$ws.Range("C122").Value = 263
